$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 421777.6
$ws.Range("J6").Value = 6375.125
$ws.Range("L6").Value = 19125.375
$ws.Range("N6").Value = -19349.375

$ws.Range("H9").Value = 171.73334
$ws.Range("I9").Value = 171.22223
$ws.Range("J9").Value = 172.5
$ws.Range("K9").Value = 171.22223
$ws.Range("L9").Value = 172.5
$ws.Range("M9").Value = -2.222229999999996
$ws.Range("N9").Value = -510.5

$ws.Range("H12").Value = 222.17647
$ws.Range("I12").Value = 216.21428
$ws.Range("J12").Value = 250
$ws.Range("K12").Value = 216.21428
$ws.Range("L12").Value = 250
$ws.Range("M12").Value = -46.21428
$ws.Range("N12").Value = -590

$ws.Range("H21").Value = 38173
$ws.Range("I21").Value = 31678.666
$ws.Range("J21").Value = 40337.777
$ws.Range("K21").Value = 31678.666
$ws.Range("L21").Value = 40337.777
$ws.Range("M21").Value = -31210.666
$ws.Range("N21").Value = -41273.777

$ws.Range("H23").Value = 38173
$ws.Range("I23").Value = 31678.666
$ws.Range("J23").Value = 40337.777
$ws.Range("K23").Value = 31678.666
$ws.Range("L23").Value = 40337.777
$ws.Range("M23").Value = -31444.666
$ws.Range("N23").Value = -40805.777

$ws.Range("H29").Value = 3130.3333
$ws.Range("I29").Value = 2299.8
$ws.Range("J29").Value = 3389.875
$ws.Range("K29").Value = 6899.400000000001
$ws.Range("L29").Value = 10169.625
$ws.Range("M29").Value = -6618.400000000001
$ws.Range("N29").Value = -10731.625

$ws.Range("H38").Value = 625415.6
$ws.Range("I38").Value = 909219.4399999999
$ws.Range("J38").Value = 1047.2
$ws.Range("K38").Value = 2727658.32
$ws.Range("L38").Value = 3141.6
$ws.Range("M38").Value = -2727286.32
$ws.Range("N38").Value = -3885.6

$ws.Range("H54").Value = 25325.334
$ws.Range("I54").Value = 15538
$ws.Range("K54").Value = 15538
$ws.Range("M54").Value = -15052

$ws.Range("H58").Value = 6445.375
$ws.Range("I58").Value = 262.5
$ws.Range("J58").Value = 10155.1
$ws.Range("K58").Value = 787.5
$ws.Range("L58").Value = 30465.3
$ws.Range("M58").Value = -637.5
$ws.Range("N58").Value = -30765.3

$ws.Range("H61").Value = 1279
$ws.Range("I61").Value = 131.66667
$ws.Range("J61").Value = 3000
$ws.Range("K61").Value = 395.00001
$ws.Range("L61").Value = 9000
$ws.Range("M61").Value = -223.00001
$ws.Range("N61").Value = -9344

$ws.Range("H81").Value = 43875
$ws.Range("J81").Value = 48500
$ws.Range("L81").Value = 48500
$ws.Range("N81").Value = -50496

$ws.Range("H84").Value = 43875
$ws.Range("J84").Value = 48500
$ws.Range("L84").Value = 145500
$ws.Range("N84").Value = -155484

$ws.Range("H98").Value = 678.92
$ws.Range("I98").Value = 719.4091
$ws.Range("J98").Value = 382
$ws.Range("K98").Value = 719.4091
$ws.Range("L98").Value = 382
$ws.Range("M98").Value = 778.5909
$ws.Range("N98").Value = -3378

$ws.Range("H122").Value = 678.92
$ws.Range("I122").Value = 719.4091
$ws.Range("J122").Value = 382
$ws.Range("K122").Value = 2158.2273
$ws.Range("L122").Value = 1146
$ws.Range("M122").Value = 291.7727
$ws.Range("N122").Value = -6046

$ws.Range("H137").Value = 25570.635
$ws.Range("I137").Value = 67531.47
$ws.Range("J137").Value = 1362.4615
$ws.Range("K137").Value = 202594.41
$ws.Range("L137").Value = 4087.3845
$ws.Range("M137").Value = -200044.41
$ws.Range("N137").Value = -9187.3845

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 4977.6665
$ws.Range("I5").Value = 6459.8
$ws.Range("J5").Value = 3125
$ws.Range("K5").Value = 6459.8
$ws.Range("L5").Value = 3125
$ws.Range("M5").Value = -6347.8
$ws.Range("N5").Value = -3349

$ws.Range("H24").Value = 46396.332
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 46396.332
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 46396.332
$ws.Range("M24").ClearContents()
$ws.Range("N24").Value = -47144.332

$ws.Range("H98").Value = 20355
$ws.Range("J98").Value = 20355
$ws.Range("L98").Value = 20355
$ws.Range("N98").Value = -26345

$ws.Range("H100").Value = 46396.332
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 46396.332
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 46396.332
$ws.Range("M100").ClearContents()
$ws.Range("N100").Value = -48560.332

$ws.Range("H110").Value = 511
$ws.Range("I110").Value = 511
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 511
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = 1534
$ws.Range("N110").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 4977.6665
$ws.Range("I4").Value = 6459.8
$ws.Range("J4").Value = 3125
$ws.Range("K4").Value = 6459.8
$ws.Range("L4").Value = 3125
$ws.Range("M4").Value = -6344.8
$ws.Range("N4").Value = -3355

$ws.Range("H95").Value = 24833.334
$ws.Range("J95").Value = 24833.334
$ws.Range("L95").Value = 24833.334
$ws.Range("N95").Value = -30325.334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H106").Value = 44960
$ws.Range("J106").Value = 44960
$ws.Range("L106").Value = 44960
$ws.Range("N106").Value = -47484

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H60").Value = 414.6
$ws.Range("I60").Value = 414.6
$ws.Range("K60").Value = 1243.8
$ws.Range("M60").Value = -992.8000000000002

$ws.Range("H76").Value = 10453.16
$ws.Range("I76").Value = 3499.5
$ws.Range("J76").Value = 11057.826
$ws.Range("K76").Value = 10498.5
$ws.Range("L76").Value = 33173.478
$ws.Range("M76").Value = -10115.5
$ws.Range("N76").Value = -33939.478

$ws.Range("H79").Value = 10453.16
$ws.Range("I79").Value = 3499.5
$ws.Range("J79").Value = 11057.826
$ws.Range("K79").Value = 10498.5
$ws.Range("L79").Value = 33173.478
$ws.Range("M79").Value = -9172.5
$ws.Range("N79").Value = -35825.478

$ws.Range("H101").Value = 4328.143
$ws.Range("I101").Value = 350
$ws.Range("J101").Value = 5919.4
$ws.Range("K101").Value = 1050
$ws.Range("L101").Value = 17758.2
$ws.Range("M101").Value = 1384
$ws.Range("N101").Value = -22626.2

$ws.Range("H103").Value = 2399.75
$ws.Range("I103").Value = 979.4286
$ws.Range("J103").Value = 3504.4443
$ws.Range("K103").Value = 2938.2858
$ws.Range("L103").Value = 10513.3329
$ws.Range("M103").Value = -2059.2858
$ws.Range("N103").Value = -12271.3329

$ws.Range("H111").Value = 2042.7
$ws.Range("I111").Value = 505.4
$ws.Range("J111").Value = 3580
$ws.Range("K111").Value = 1516.2
$ws.Range("L111").Value = 10740
$ws.Range("M111").Value = 1550.8
$ws.Range("N111").Value = -16874

$ws.Range("H113").Value = 670.6
$ws.Range("I113").Value = 800
$ws.Range("J113").Value = 656.2222
$ws.Range("K113").Value = 2400
$ws.Range("L113").Value = 1968.6666
$ws.Range("M113").Value = -230
$ws.Range("N113").Value = -6308.6666

$ws.Range("H137").Value = 3097466
$ws.Range("I137").Value = 1261.8182
$ws.Range("J137").Value = 4459795.5
$ws.Range("K137").Value = 3785.4546
$ws.Range("L137").Value = 13379386.5
$ws.Range("M137").Value = 1314.5454
$ws.Range("N137").Value = -13389586.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2615.3845
$ws.Range("I80").Value = 3000
$ws.Range("J80").Value = 2444.4443
$ws.Range("K80").Value = 3000
$ws.Range("L80").Value = 2444.4443
$ws.Range("M80").Value = -2002
$ws.Range("N80").Value = -4440.4443

$ws.Range("H83").Value = 2615.3845
$ws.Range("I83").Value = 3000
$ws.Range("J83").Value = 2444.4443
$ws.Range("K83").Value = 15000
$ws.Range("L83").Value = 12222.2215
$ws.Range("M83").Value = -10008
$ws.Range("N83").Value = -22206.2215

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H97").Value = 23586
$ws.Range("J97").Value = 23586
$ws.Range("L97").Value = 23586
$ws.Range("N97").Value = -25568

$ws.Range("H101").Value = 5580
$ws.Range("J101").Value = 5580
$ws.Range("L101").Value = 5580
$ws.Range("N101").Value = -12070

$ws.Range("H105").Value = 38000
$ws.Range("J105").Value = 38000
$ws.Range("L105").Value = 38000
$ws.Range("N105").Value = -44988

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 14333.333
$ws.Range("J63").Value = 14333.333
$ws.Range("L63").Value = 14333.333
$ws.Range("N63").Value = -15581.333

$ws.Range("H66").Value = 14333.333
$ws.Range("J66").Value = 14333.333
$ws.Range("L66").Value = 42999.999
$ws.Range("N66").Value = -49239.999

$ws.Range("H81").Value = 1326.9166
$ws.Range("I81").Value = 637.5
$ws.Range("J81").Value = 1556.7222
$ws.Range("K81").Value = 1275
$ws.Range("L81").Value = 3113.4444
$ws.Range("M81").Value = -214
$ws.Range("N81").Value = -5235.4444

$ws.Range("H84").Value = 1326.9166
$ws.Range("I84").Value = 637.5
$ws.Range("J84").Value = 1556.7222
$ws.Range("K84").Value = 6375
$ws.Range("L84").Value = 15567.222
$ws.Range("M84").Value = -1071
$ws.Range("N84").Value = -26175.222
